$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new header cell "Result For" in A1, matching B1's bold/bordered style
$ws.Range("A1").Value = "Result For"
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)  # xlPasteFormats

# Change the status value from "passed" to "PASS"
$ws.Range("C11").Value = "PASS"

# Move the active selection to C11
$ws.Range("C11").Select()

# Adjust column B width to match new content width
$ws.Columns.Item(2).ColumnWidth = 52.3
